$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column D: "preco" ---
$ws.Range("D1").Value = "preco"
# Copy header style (bold, centered, bordered) from C1 onto D1
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the preco values for existing rows 2-9
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 44
$ws.Range("D6").Value = 200
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 0

# --- Update existing quantidade_disponivel values ---
$ws.Range("B2").Value = 0
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 3000

# --- New rows ---
$ws.Range("A10").Value = "tomate seco timy sache"
$ws.Range("B10").Value = 3000
$ws.Range("C10").Value = "g"
$ws.Range("D10").Value = 50

$ws.Range("A11").Value = "azeitona"
$ws.Range("B11").Value = 100
$ws.Range("C11").Value = "Un"
$ws.Range("D11").Value = 12
